# Dragon & Phoenix review doc:
#  1. Remove the old "Meta description: ..." paragraph that sat right
#     under the Heading1 title.
#  2. Turn the trailing DALLE image-prompt paragraph into two paragraphs:
#     a bold "Play Dragon & Phoenix Slot for Free - Betsoft 2019 Game"
#     line, followed by the (still italic) meta-description copy that
#     used to live at the top of the document.

$d = $word.ActiveDocument

# --- 1. Delete the "Meta description: ..." paragraph -----------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description:*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -ge 1) {
    $d.Paragraphs.Item($metaIndex).Range.Delete()
}

# --- 2. Locate the DALLE image-prompt paragraph -----------------------
$dalleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "DALLE, please create a feature image*") {
        $dalleIndex = $i
        break
    }
}

if ($dalleIndex -ge 1) {
    # Insert a new paragraph right before the DALLE paragraph by typing
    # the title text + paragraph mark immediately after the end of the
    # PRECEDING paragraph. Doing it this way (rather than splitting the
    # DALLE paragraph itself) means the new text does not inherit the
    # DALLE paragraph's italic run formatting.
    $prevPara = $d.Paragraphs.Item($dalleIndex - 1)
    $insertPoint = $d.Range($prevPara.Range.End, $prevPara.Range.End)
    $insertPoint.InsertAfter("Play Dragon & Phoenix Slot for Free - Betsoft 2019 Game`r")

    # The new title paragraph now occupies the DALLE paragraph's old
    # index; bold its text (leave the paragraph mark itself alone).
    $titlePara = $d.Paragraphs.Item($dalleIndex)
    $titleTextRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
    $titleTextRange.Font.Bold = $true

    # The DALLE paragraph shifted down by one; swap its text for the new
    # meta-description copy while keeping its existing italic run.
    $dallePara = $d.Paragraphs.Item($dalleIndex + 1)
    $dalleTextRange = $d.Range($dallePara.Range.Start, $dallePara.Range.End - 1)
    $dalleTextRange.Text = "Read our Dragon and Phoenix slot review and play the game for free. Features, pros, and cons, and potential for big wins up to 8098x total bet."
}
